$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create the new "2022-Q4" sheet by duplicating "2022-Q3" (so it inherits
#    the same layout / header styling) and placing it immediately before it.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3, $null)
$q4 = $wb.ActiveSheet
$q4.Name = "2022-Q4"

# Force text storage ("@") for the numeric-looking text columns (fund code,
# fund size, position, etc.) on the rows we are about to fill in, otherwise
# values such as "004685" would be silently re-interpreted as numbers.
$q4.Range("B2:G4").NumberFormat = "@"

# Row 2 - 金元顺安元启灵活配置混合
$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "004685"
$q4.Range("C2").Value = "金元顺安元启灵活配置混合"
$q4.Range("D2").Value = "15.29"
$q4.Range("E2").Value = "76.11"
$q4.Range("F2").Value = "1.10"
$q4.Range("G2").Value = "0.1682"
$q4.Range("H2").Value = 1

# Row 3 - 大成中证360互联网+大数据100指数A
$q4.Range("A2").Copy()
$q4.Range("A3").PasteSpecial(-4122)
$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "002236"
$q4.Range("C3").Value = "大成中证360互联网+大数据100指数A"
$q4.Range("D3").Value = "1.15"
$q4.Range("E3").Value = "92.50"
$q4.Range("F3").Value = "1.02"
$q4.Range("G3").Value = "0.0117"
$q4.Range("H3").Value = 3

# Row 4 - 大成中证360互联网+大数据100指数C
$q4.Range("A2").Copy()
$q4.Range("A4").PasteSpecial(-4122)
$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "003359"
$q4.Range("C4").Value = "大成中证360互联网+大数据100指数C"
$q4.Range("D4").Value = "1.12"
$q4.Range("E4").Value = "92.50"
$q4.Range("F4").Value = "1.02"
$q4.Range("G4").Value = "0.0114"
$q4.Range("H4").Value = 3

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: push the existing two rows down by
#    one and insert the new "2022-Q4" totals in row 2.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$oldQ2Date = $total.Range("B3").Value()
$oldQ2Count = $total.Range("C3").Value()
$oldQ2Aum = $total.Range("D3").Value()

$oldQ3Date = $total.Range("B2").Value()
$oldQ3Count = $total.Range("C2").Value()
$oldQ3Aum = $total.Range("D2").Value()

# Row 4 <- old row 3 (2022-Q2), carry the A-column index style along.
$total.Range("A2").Copy()
$total.Range("A4").PasteSpecial(-4122)
$total.Range("A4").Value = 2
$total.Range("B4").Value = $oldQ2Date
$total.Range("C4").Value = $oldQ2Count
$total.Range("D4").Value = $oldQ2Aum

# Row 3 <- old row 2 (2022-Q3)
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)
$total.Range("A3").Value = 1
$total.Range("B3").Value = $oldQ3Date
$total.Range("C3").Value = $oldQ3Count
$total.Range("D3").Value = $oldQ3Aum

# Row 2 <- new 2022-Q4 totals
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.19
